# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts.
#
# This updates the DAMSLTag (column I) and DialogAct (column J) values
# for a set of rows in the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (row 1 headers):
#   I = DAMSLTag
#   J = DialogAct
$colDamsl = 9
$colDialogAct = 10

# Map of spreadsheet row -> new (DAMSLTag, DialogAct) values.
$updates = @{
    3   = @("sv", "Statement-opinion")
    7   = @("sv", "Statement-opinion")
    9   = @("sv", "Statement-opinion")
    12  = @("sd", "Statement-non-opinion")
    14  = @("aa", "Agree/Accept")
    34  = @("b", "Acknowledge (Backchannel)")
    35  = @("sd", "Statement-non-opinion")
    40  = @("sv", "Statement-opinion")
    55  = @("sv", "Statement-opinion")
    56  = @("sv", "Statement-opinion")
    58  = @("b", "Acknowledge (Backchannel)")
    60  = @("sv", "Statement-opinion")
    68  = @("b", "Acknowledge (Backchannel)")
    87  = @("aa", "Agree/Accept")
    93  = @("sd", "Statement-non-opinion")
    105 = @("sd", "Statement-non-opinion")
    109 = @("sd", "Statement-non-opinion")
    110 = @("sd", "Statement-non-opinion")
    119 = @("sd", "Statement-non-opinion")
    122 = @("sd", "Statement-non-opinion")
    137 = @("sv", "Statement-opinion")
    152 = @("sv", "Statement-opinion")
    156 = @("sv", "Statement-opinion")
    158 = @("sv", "Statement-opinion")
    201 = @("%", "Uninterpretable")
    203 = @("sd", "Statement-non-opinion")
    210 = @("sd", "Statement-non-opinion")
    217 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, $colDamsl).Value = $values[0]
    $ws.Cells.Item($row, $colDialogAct).Value = $values[1]
}
